# Auto-generated edit script: apply updated Leve profit data values
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per scheduled runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 111
$ws.Range("I5").Value = 99.875
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 99.875
$ws.Range("L5").Value = 200
$ws.Range("M5").Value = 15.125
$ws.Range("N5").Value = -430
$ws.Range("H12").Value = 497.5
$ws.Range("I12").Value = 500
$ws.Range("J12").Value = 495
$ws.Range("K12").Value = 500
$ws.Range("L12").Value = 495
$ws.Range("M12").Value = -330
$ws.Range("N12").Value = -835
$ws.Range("H18").Value = 1450
$ws.Range("I18").Value = 1450
$ws.Range("K18").Value = 1450
$ws.Range("M18").Value = -1166
$ws.Range("H51").Value = 9000
$ws.Range("I51").Value = 9000
$ws.Range("K51").Value = 9000
$ws.Range("M51").Value = -8516
$ws.Range("H62").Value = 8762.538
$ws.Range("J62").Value = 10665.223
$ws.Range("L62").Value = 10665.223
$ws.Range("N62").Value = -11913.223
$ws.Range("H65").Value = 8762.538
$ws.Range("J65").Value = 10665.223
$ws.Range("L65").Value = 53326.115
$ws.Range("N65").Value = -59566.115
$ws.Range("H76").Value = 4700.4
$ws.Range("I76").Value = 3667.3333
$ws.Range("K76").Value = 3667.3333
$ws.Range("M76").Value = -3352.3333
$ws.Range("H79").Value = 4700.4
$ws.Range("I79").Value = 3667.3333
$ws.Range("K79").Value = 3667.3333
$ws.Range("M79").Value = -2575.3333
$ws.Range("H86").Value = 2100
$ws.Range("I86").Value = 4000
$ws.Range("J86").Value = 200
$ws.Range("K86").Value = 4000
$ws.Range("L86").Value = 200
$ws.Range("M86").Value = -2877
$ws.Range("N86").Value = -2446
$ws.Range("H89").Value = 2100
$ws.Range("I89").Value = 4000
$ws.Range("J89").Value = 200
$ws.Range("K89").Value = 20000
$ws.Range("L89").Value = 1000
$ws.Range("M89").Value = -14384
$ws.Range("N89").Value = -12232
$ws.Range("H95").Value = 28655.5
$ws.Range("J95").Value = 28655.5
$ws.Range("L95").Value = 28655.5
$ws.Range("N95").Value = -34147.5
$ws.Range("H106").Value = 6500.3335
$ws.Range("I106").Value = 6100.4
$ws.Range("K106").Value = 6100.4
$ws.Range("M106").Value = -5469.4
$ws.Range("H107").Value = 128.8
$ws.Range("I107").Value = 128.8
$ws.Range("K107").Value = 128.8
$ws.Range("M107").Value = 1791.2
$ws.Range("H112").Value = 942.5
$ws.Range("I112").Value = 1295
$ws.Range("K112").Value = 3885
$ws.Range("M112").Value = -2777
$ws.Range("H116").Value = 3340
$ws.Range("I116").Value = 1900
$ws.Range("K116").Value = 1900
$ws.Range("M116").Value = 1542
$ws.Range("H132").Value = 13923.105
$ws.Range("I132").Value = 13271.1875
$ws.Range("K132").Value = 39813.5625
$ws.Range("M132").Value = -37283.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 111
$ws.Range("I4").Value = 111
$ws.Range("K4").Value = 111
$ws.Range("M4").Value = 5
$ws.Range("H110").Value = 3703.3635
$ws.Range("J110").Value = 4171
$ws.Range("L110").Value = 4171
$ws.Range("N110").Value = -8261

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3549.3572
$ws.Range("I20").Value = 4208.6
$ws.Range("K20").Value = 4208.6
$ws.Range("M20").Value = -3961.6
$ws.Range("H86").Value = 3505.4736
$ws.Range("I86").Value = 1459.8889
$ws.Range("J86").Value = 5346.5
$ws.Range("K86").Value = 1459.8889
$ws.Range("L86").Value = 5346.5
$ws.Range("M86").Value = -336.8888999999999
$ws.Range("N86").Value = -7592.5
$ws.Range("H89").Value = 3505.4736
$ws.Range("I89").Value = 1459.8889
$ws.Range("J89").Value = 5346.5
$ws.Range("K89").Value = 7299.4445
$ws.Range("L89").Value = 26732.5
$ws.Range("M89").Value = -1683.4445
$ws.Range("N89").Value = -37964.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1561.7368
$ws.Range("I22").Value = 898.63635
$ws.Range("J22").Value = 2473.5
$ws.Range("K22").Value = 898.63635
$ws.Range("L22").Value = 2473.5
$ws.Range("M22").Value = -548.63635
$ws.Range("N22").Value = -3173.5
$ws.Range("H134").Value = 2007.6
$ws.Range("I134").Value = 2007.6
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 6022.799999999999
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -3487.799999999999
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 386
$ws.Range("I14").Value = 386
$ws.Range("K14").Value = 1158
$ws.Range("M14").Value = -985
$ws.Range("H56").Value = 6666.6665
$ws.Range("I56").Value = 6666.6665
$ws.Range("K56").Value = 6666.6665
$ws.Range("M56").Value = -6136.6665
$ws.Range("H107").Value = 729
$ws.Range("J107").Value = 729
$ws.Range("L107").Value = 2187
$ws.Range("N107").Value = -6027
$ws.Range("H108").Value = 424.5
$ws.Range("I108").Value = 424.5
$ws.Range("K108").Value = 1273.5
$ws.Range("M108").Value = 1606.5
$ws.Range("H109").Value = 772179.6
$ws.Range("I109").Value = 2001487.2
$ws.Range("J109").Value = 3862.375
$ws.Range("K109").Value = 6004461.6
$ws.Range("L109").Value = 11587.125
$ws.Range("M109").Value = -6003421.6
$ws.Range("N109").Value = -13667.125
$ws.Range("H140").Value = 4740
$ws.Range("I140").Value = 3200
$ws.Range("K140").Value = 9600
$ws.Range("M140").Value = -4420

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 4752.25
$ws.Range("I55").Value = 4904.8335
$ws.Range("K55").Value = 4904.8335
$ws.Range("M55").Value = -4577.8335
$ws.Range("H70").Value = 5590.25
$ws.Range("I70").Value = 2297.25
$ws.Range("J70").Value = 12176.25
$ws.Range("K70").Value = 2297.25
$ws.Range("L70").Value = 12176.25
$ws.Range("M70").Value = -2027.25
$ws.Range("N70").Value = -12716.25
$ws.Range("H73").Value = 5590.25
$ws.Range("I73").Value = 2297.25
$ws.Range("J73").Value = 12176.25
$ws.Range("K73").Value = 2297.25
$ws.Range("L73").Value = 12176.25
$ws.Range("M73").Value = -1361.25
$ws.Range("N73").Value = -14048.25
$ws.Range("H122").Value = 502565.3
$ws.Range("I122").Value = 627331.9
$ws.Range("K122").Value = 1881995.7
$ws.Range("M122").Value = -1879545.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4887.074
$ws.Range("J46").Value = 6008.3335
$ws.Range("L46").Value = 6008.3335
$ws.Range("N46").Value = -6384.3335
$ws.Range("H100").Value = 7481.2
$ws.Range("I100").Value = 2456
$ws.Range("J100").Value = 8737.5
$ws.Range("K100").Value = 2456
$ws.Range("L100").Value = 8737.5
$ws.Range("M100").Value = -1915
$ws.Range("N100").Value = -9819.5
$ws.Range("H106").Value = 20396.75
$ws.Range("J106").Value = 20396.75
$ws.Range("L106").Value = 20396.75
$ws.Range("N106").Value = -22920.75
$ws.Range("H122").Value = 4894.5
$ws.Range("I122").Value = 4868.375
$ws.Range("K122").Value = 14605.125
$ws.Range("M122").Value = -12155.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6872.727
$ws.Range("J62").Value = 8885.857
$ws.Range("L62").Value = 8885.857
$ws.Range("N62").Value = -10133.857
$ws.Range("H65").Value = 6872.727
$ws.Range("J65").Value = 8885.857
$ws.Range("L65").Value = 44429.285
$ws.Range("N65").Value = -50669.285
$ws.Range("H81").Value = 1114.9
$ws.Range("I81").Value = 731.125
$ws.Range("J81").Value = 2650
$ws.Range("K81").Value = 1462.25
$ws.Range("L81").Value = 5300
$ws.Range("M81").Value = -401.25
$ws.Range("N81").Value = -7422
$ws.Range("H84").Value = 1114.9
$ws.Range("I84").Value = 731.125
$ws.Range("J84").Value = 2650
$ws.Range("K84").Value = 7311.25
$ws.Range("L84").Value = 26500
$ws.Range("M84").Value = -2007.25
$ws.Range("N84").Value = -37108
$ws.Range("H100").Value = 3250.6667
$ws.Range("I100").Value = 3250.6667
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 6501.3334
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -5960.3334
$ws.Range("N100").ClearContents()
$ws.Range("H107").Value = 569.4
$ws.Range("I107").Value = 436.75
$ws.Range("K107").Value = 1310.25
$ws.Range("M107").Value = 609.75

Write-Host "Applied Leve profit data updates."
